$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Customer" sheet right after "Driver"
# ---------------------------------------------------------------------------
$driver = $wb.Worksheets.Item("Driver")
$customer = $wb.Worksheets.Add($null, $driver)
$customer.Name = "Customer"

# Header row
$customer.Range("A1").Value = "CustomerId"
$customer.Range("B1").Value = "FirstName"
$customer.Range("C1").Value = "LastName"
$customer.Range("D1").Value = "Birthdate"

# Numeric Id column
$customer.Range("A2").Value = 1
$customer.Range("A3").Value = 2
$customer.Range("A4").Value = 3
$customer.Range("A5").Value = 4
$customer.Range("A6").Value = 5
$customer.Range("A7").Value = 6
$customer.Range("A8").Value = 7
$customer.Range("A9").Value = 8
$customer.Range("A10").Value = 9
$customer.Range("A11").Value = 10
$customer.Range("A12").Value = 11
$customer.Range("A13").Value = 12
$customer.Range("A14").Value = 13
$customer.Range("A15").Value = 14
$customer.Range("A16").Value = 15
$customer.Range("A17").Value = 16
$customer.Range("A18").Value = 17
$customer.Range("A19").Value = 18
$customer.Range("A20").Value = 19
$customer.Range("A21").Value = 20

# Customer rows are the Driver rows copied in this order: 11-21, 6-10, 2-5
$driver.Range("B11:D21").Copy()
$customer.Range("B2").PasteSpecial(-4163)

$driver.Range("B6:D10").Copy()
$customer.Range("B13").PasteSpecial(-4163)

$driver.Range("B2:D5").Copy()
$customer.Range("B18").PasteSpecial(-4163)

# Reuse the existing number/column formatting instead of inventing new styles
$driver.Range("A1").Copy()
$customer.Range("A1").PasteSpecial(-4122)

$driver.Range("D2").Copy()
$customer.Range("D2:D21").PasteSpecial(-4122)

$customer.Columns.Item(4).ColumnWidth = 16

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Rename "Routes" to "Route"
# ---------------------------------------------------------------------------
$routes = $wb.Worksheets.Item("Routes")
$routes.Name = "Route"

# ---------------------------------------------------------------------------
# 3. Reset / update the saved view state of the various sheets so it matches
#    what was left behind by the editing session
# ---------------------------------------------------------------------------
$driver.Activate()
$driver.Range("B2:D5").Select()

$truckAndDriver = $wb.Worksheets.Item("TruckAndDriver")
$truckAndDriver.Activate()
$truckAndDriver.Range("G4").Select()

$warehouse = $wb.Worksheets.Item("Warehouse")
$warehouse.Activate()
$warehouse.Range("G4").Select()

$routes.Activate()
$routes.Range("D16").Select()

# Leave "Customer" as the active / selected sheet and cell, matching the
# final state captured in the commit
$customer.Activate()
$customer.Range("F6").Select()
